# Apply the edit described by the diff:
#  1. "TP3-Log3430" becomes "TP6-Log3430", split across three runs
#     ("TP", "6", "-Log3430") with the "_GoBack" bookmark now sitting
#     between the "6" run and the "-Log3430" run.
#  2. The "_GoBack" bookmark that used to sit right after "12 avril"
#     is gone (Word only ever keeps one "_GoBack" bookmark, tracking
#     the most recent edit location, so re-adding it at the new spot
#     removes it from the old one).

$d = $word.ActiveDocument

# Locate the title text "TP3-Log3430" wherever it lives in the story.
$titleRange = $d.Content
$titleRange.Find.Execute("TP3-Log3430", $false, $false, $false, $false, `
                          $false, $true, 1, $false, "", 0) | Out-Null

$titleStart = $titleRange.Start

# Change the digit "3" -> "6" (keeps everything else untouched).
$digitRange = $d.Range($titleStart + 2, $titleStart + 3)
$digitRange.Text = "6"

# Force a run split right after "TP" (before the "6") using a
# transient bookmark, then remove the bookmark once the split exists.
$splitPoint = $d.Range($titleStart + 2, $titleStart + 2)
$d.Bookmarks.Add("__tmp_split__", $splitPoint) | Out-Null

# Re-anchor "_GoBack" right after the "6" (i.e. before "-Log3430").
# Because bookmark names are unique, this automatically removes the
# previous "_GoBack" bookmark that sat after "12 avril".
$goBackPoint = $d.Range($titleStart + 3, $titleStart + 3)
$d.Bookmarks.Add("_GoBack", $goBackPoint) | Out-Null

# Drop the helper bookmark; the run split it created remains.
$d.Bookmarks("__tmp_split__").Delete()
